# Updates the cryptos list (prices / 1h volume % / a few row-ordering swaps)
# as refreshed by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns are stored as plain text (e.g. "26.745.83" uses
# dots as thousands separators, not a number) so force text formatting before
# writing the new values - otherwise Excel would "helpfully" reinterpret a
# value like "21.15" as a numeric 21.15 and mangle others. Reset the style
# back to Normal afterwards so no stray cell formatting is left behind.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.738.77"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.537.88"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "205.55"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.244"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "21.15"
$ws.Range("E9").Value = "  -3.89%  "
$ws.Range("D10").Value = "0.0579"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").Value = "0.0854"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "1.764.51"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D13").Value = "1.561.26"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "3.66"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "0.505"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").Value = "26.699.41"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "60.99"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0₃0686"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "211.56"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").Value = "7.21"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "4.03"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").Value = "8.99"
$ws.Range("E23").Value = "  -3.95%  "
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").Value = "152.65"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "6.44"
$ws.Range("E26").Value = "  -4.38%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "14.81"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").Value = "0.0454"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("D32").Value = "3.21"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").Value = "1.355.85"
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("D34").Value = "2.91"
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "1.50"
$ws.Range("E35").Value = "  -3.38%  "
$ws.Range("D36").Value = "2.28"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").Value = "0.932"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").Value = "0.0163"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").Value = "0.520"
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "5.73"
$ws.Range("E40").Value = "  +5.22%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "0.797"
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("D43").Value = "2.19"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "1.73"
$ws.Range("E44").Value = "  -2.75%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "62.32"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("D46").Value = "1.677.49"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").Value = "2.24"
$ws.Range("E47").Value = "  -4.34%  "
$ws.Range("D48").Value = "85.78"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "0.0510"
$ws.Range("E49").Value = "  +3.06%  "
$ws.Range("D50").Value = "0.0947"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  -0.02%  "

$dataRange.Style = "Normal"

